# Update RunTestSuite and LoginTest
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# LoginTest: change Test Class in B3 from the WMC account management test
# to the patient.tests class (same as row 2) so the unused string is dropped
# from the shared string table.
$ws.Range("B3").Value = "patient.tests"

# Move the active selection on the Scenarios sheet to D4.
$ws.Range("D4").Select()
